$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused zVel/RelSpeed columns (G:H), shifting dimension to A1:F13
$ws.Range("G1:H13").Delete()

# Update header row
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "C/No"
$ws.Range("D1").Value = "Eb/No"
$ws.Range("E1").Value = "BER"
$ws.Range("F1").Value = "Range"

# Row 2
$ws.Range("B2").Value = "22 Jun 2025 10:15:29.000000000"
$ws.Range("C2").Value = -2.94901106364583
$ws.Range("D2").Value = -5.959311020285641
$ws.Range("E2").Value = 0.238197132938244
$ws.Range("F2").Value = 2754.089201487097

# Row 3
$ws.Range("B3").Value = "22 Jun 2025 10:15:39.000000000"
$ws.Range("C3").Value = -28.15247928870762
$ws.Range("D3").Value = -31.16277924534743
$ws.Range("E3").Value = 0.4843981768018749
$ws.Range("F3").Value = 2685.571758676042

# Row 4
$ws.Range("B4").Value = "22 Jun 2025 10:15:49.000000000"
$ws.Range("C4").Value = -1.229907376169047
$ws.Range("D4").Value = -4.240207332808859
$ws.Range("E4").Value = 0.1927054244717706
$ws.Range("F4").Value = 2617.106297974157

# Row 5
$ws.Range("B5").Value = "22 Jun 2025 10:15:59.000000000"
$ws.Range("C5").Value = -10.14464047948895
$ws.Range("D5").Value = -13.15494043612876
$ws.Range("E5").Value = 0.3778982457374752
$ws.Range("F5").Value = 2548.705801581525

# Row 6
$ws.Range("B6").Value = "22 Jun 2025 10:16:09.000000000"
$ws.Range("C6").Value = -10.05917224203106
$ws.Range("D6").Value = -13.06947219867087
$ws.Range("E6").Value = 0.3767298360609527
$ws.Range("F6").Value = 2480.3844497692

# Row 7
$ws.Range("B7").Value = "22 Jun 2025 10:16:19.000000000"
$ws.Range("C7").Value = -6.540831652997997
$ws.Range("D7").Value = -9.55113160963781
$ws.Range("E7").Value = 0.3188445644709756
$ws.Range("F7").Value = 2412.157784091685

# Row 8
$ws.Range("B8").Value = "22 Jun 2025 10:16:29.000000000"
$ws.Range("C8").Value = -19.53143848642713
$ws.Range("D8").Value = -22.54173844306695
$ws.Range("E8").Value = 0.4579726067351606
$ws.Range("F8").Value = 2344.042903608837

# Row 9
$ws.Range("B9").Value = "22 Jun 2025 10:16:39.000000000"
$ws.Range("C9").Value = -1.744837658402866
$ws.Range("D9").Value = -4.755137615042678
$ws.Range("E9").Value = 0.2066760055034322
$ws.Range("F9").Value = 2276.058681724883

# Row 10
$ws.Range("B10").Value = "22 Jun 2025 10:16:49.000000000"
$ws.Range("C10").Value = -7.388618359481662
$ws.Range("D10").Value = -10.39891831612147
$ws.Range("E10").Value = 0.3346390959397672
$ws.Range("F10").Value = 2208.226032162438

# Row 11
$ws.Range("B11").Value = "22 Jun 2025 10:16:59.000000000"
$ws.Range("C11").Value = -0.8371421059926083
$ws.Range("D11").Value = -3.847442062632421
$ws.Range("E11").Value = 0.1819076049461735
$ws.Range("F11").Value = 2140.568221581015

# Row 12
$ws.Range("B12").Value = "22 Jun 2025 10:17:09.000000000"
$ws.Range("C12").Value = -2.838639365189925
$ws.Range("D12").Value = -5.848939321829737
$ws.Range("E12").Value = 0.2353869473907544
$ws.Range("F12").Value = 2073.111225564506

# Row 13
$ws.Range("B13").Value = "22 Jun 2025 10:17:19.000000000"
$ws.Range("C13").Value = -3.892056142929642
$ws.Range("D13").Value = -6.902356099569453
$ws.Range("E13").Value = 0.2614610754993032
$ws.Range("F13").Value = 2005.884170588258
